$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the values that used to live in row 3 (D2),
# plus freshly-simulated results for H2/I2/J2.
$ws.Range("D2").Value = 1000000
$ws.Range("H2").Value = 28917147085.4152
$ws.Range("I2").Value = 6791.110624416735
$ws.Range("J2").Value = 21770108661.10713

# Row 3 (the old duplicate/second group run) is no longer needed -
# remove it entirely so the data block shrinks from A1:J3 to A1:J2.
$ws.Rows("3").Delete()
